# This script merges several runs that hold adjacent pieces of the same
# sentence/phrase (e.g. "NOAA" + " " + "quarto" + " " + "book") into a
# single run with the full text, mirroring what happens when Word's
# renderer/round-trip collapses runs that carry identical formatting.
#
# We do this with Find/Replace: searching for the full (already
# contiguous, when concatenated across runs) text and replacing it with
# the very same text forces Word to rewrite the matched range as one run.

$d = $word.ActiveDocument

function Merge-Text($needle) {
    $result = $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $needle"
    }
}

Merge-Text("NOAA quarto book")
Merge-Text("Jane Doe")
Merge-Text("Eva Nováková")
Merge-Text("Matti Meikäläinen")
Merge-Text("Table of contents")
Merge-Text([string]::Concat([char]0x201C, "use template", [char]0x201D))
Merge-Text([string]::Concat([char]0x201C, "Version Control", [char]0x201D))
Merge-Text([string]::Concat([char]0x201C, "Build", [char]0x201D))
Merge-Text([string]::Concat([char]0x201C, "Render website", [char]0x201D))
Merge-Text([string]::Concat([char]0x201C, "preview in browser", [char]0x201D))
Merge-Text([string]::Concat([char]0x201C, "Show output preview in: Viewer panel", [char]0x201D))

Write-Host "done"
